$d = $word.ActiveDocument

# The last paragraph currently reads:
#   "la methode notifyObserves (Object arg)"
# It needs its trailing run split (spell-check re-ran over "arg"), and a brand
# new paragraph describing another bug report needs to be appended right after it.
# Rebuild the paragraph + append the new one as literal OOXML so the run /
# w:proofErr boundaries come out exactly as Word's proofing pass would produce.

$lastPara = $d.Paragraphs.Last
$lastPara.Range.Delete()

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00D45D62" w:rsidRDefault="0046258D" w:rsidP="00D45D62"><w:pPr><w:pStyle w:val="Sansinterligne"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>methode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>notifyObserves</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (Object </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>arg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Sansinterligne"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">j’ai un </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pb</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> aussi a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cahque</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> fois que clique sur le bouton il m’ouvre une nouvelle </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fenetre</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> d’affichage a chaque clique</w:t></w:r></w:p>'

$insertionPoint = $d.Content
$insertionPoint.Collapse(0)
$null = $insertionPoint.InsertXML($newXml)
